# AutomationDataRetanqueo.xlsx — "ajuste pesp y firma documentos codigo ado digicredito"
#
# The edit touches the "RetanqueoDigiCreditoCCS" sheet (row 2, the sample
# data row used by the automated test), updating several fields to a new
# test scenario (new pagaduria, new disbursement date, new signature date,
# new saneamiento amount, new contact data, new credit line, etc.), plus
# moving the active selection on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RetanqueoDigiCreditoCCS")
$ws.Select()

# --- Row 2 data updates -----------------------------------------------
# Written in the same order the original authoring session used (matters
# for shared-string insertion order on save), ending with the signature
# date (P2) last.

# Pagaduria: "COLFONDOS" -> "P.A COLPENSIONES"
$ws.Range("A2").Value = '"P.A COLPENSIONES"'

# DiasHabilesIntereses: "50" -> "120"
$ws.Range("F2").Value = '"120"'

# fechaDesembolso: "07/03/2022" -> "14/03/2022"
$ws.Range("K2").Value = '"14/03/2022"'

# vlrCompasSaneamientos: "300000" -> "450000"
$ws.Range("T2").Value = '"450000"'

# Correo: "daabogadog@hotmail.com" -> "daabogadog@gmail.com" (trailing space kept)
$ws.Range("W2").Value = '"daabogadog@gmail.com" '

# Celular: "3142739038" -> "3115139037"
$ws.Range("X2").Value = '"3115139037"'

# AnoAnalisis: "2021" -> "2022"
$ws.Range("AL2").Value = '"2022"'

# Cartera1: "300000" -> "0"
$ws.Range("AR2").Value = '"0"'

# lineaCredito: "Retanqueo libre inversión" -> "Retanqueo compra de cartera"
$ws.Range("AX2").Value = '"Retanqueo compra de cartera"'

# tomarSeguroAP: "No" -> "Si"
$ws.Range("CW2").Value = '"Si"'

# NumRadicacion (codigo ADO): "92845" -> "92862"
$ws.Range("DM2").Value = '"92862"'

# fechaActual (firma documentos): "09/03/2022" -> "18/03/2022"
$ws.Range("P2").Value = '"18/03/2022"'

# --- View / selection ---------------------------------------------------
# Move the visible window / active cell to P6 (scrolled so column I is
# the left-most visible column), mirroring the saved selection state.
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("P6").Select()
